{"js": "// Rewrite the M2Doc template field `{ m:null.resize(null) }`, which is\n// currently stored as a live Word field (fldChar begin/instrText/fldChar\n// end), into plain literal text runs spelling out the same field syntax:\n// \"{\", \"m\", \":\", \"null\", \".resize(\", \"null\", \")\", \"}\" \u2014 one run per\n// original instrText run, keeping the orange theme color on the\n// \"null.resize(null)\" portion.\n\nconst body = context.document.body;\n\n// Locate the field whose code matches the M2Doc \"m:null.resize(null)\" query,\n// and the paragraph that owns it.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet targetField = null;\nlet targetParagraph = null;\n\nfor (const p of paragraphs.items) {\n    const pFields = p.getRange(\"Content\").fields;\n    pFields.load(\"items/code\");\n    await context.sync();\n\n    const match = pFields.items.find(\n        f => f.code && f.code.indexOf(\"null.resize(null)\") !== -1\n    );\n    if (match) {\n        targetField = match;\n        targetParagraph = p;\n        break;\n    }\n}\n\nif (targetField && targetParagraph) {\n    // Remove the live field (fldChar begin/instrText/fldChar end) first,\n    // leaving the (now empty) host paragraph untouched otherwise.\n    targetField.delete();\n    await context.sync();\n\n    // Re-fetch the paragraph's content range now that the field is gone.\n    const range = targetParagraph.getRange(\"Content\");\n\n    // Flat-OPC fragment with literal text runs replacing the field. The\n    // paragraph's own identity attributes are carried over explicitly so\n    // that even though InsertXML/Replace swaps in a brand-new <w:p>, the\n    // paragraph keeps the same identity as before.\n    const ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n        '<w:body>' +\n        '<w:p w14:paraId=\"2ACC4236\" w14:textId=\"4215FF0E\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n        '<w:r><w:t>{</w:t></w:r>' +\n        '<w:r><w:t>m</w:t></w:r>' +\n        '<w:r><w:t>:</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.resize(</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>)</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>';\n\n    range.insertOoxml(ooxml, Word.InsertLocation.replace);\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the field that contains the \"m:null.resize(null)\" template code and\n# rewrite it from a live Word field (fldChar/instrText) into plain literal\n# text runs: \"{\", \"m\", \":\", \"null\", \".resize(\", \"null\", \")\", \"}\" \u2014 splitting\n# the run boundaries exactly like the original instrText runs, and keeping\n# the orange theme color on the \"null.resize(null)\" part.\n$target = $null\nforeach ($f in $d.Fields) {\n    if ($f.Code.Text -match \"m:null\\.resize\\(null\\)\") {\n        $target = $f\n    }\n}\n\nif ($target -ne $null) {\n    $fieldStart = $target.Code.Start\n\n    $p = $null\n    foreach ($cand in $d.Paragraphs) {\n        if ($cand.Range.Start -le $fieldStart -and $cand.Range.End -ge $fieldStart) {\n            $p = $cand\n        }\n    }\n\n    $target.Delete()\n\n    $r = $p.Range\n\n    $xmlFrag = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r><w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.resize(</w:t></w:r><w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r><w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>)</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n    # Insert into the (now empty) paragraph, excluding the trailing paragraph\n    # mark, so the paragraph's own properties/attributes are preserved and\n    # only its runs are replaced.\n    $rr = $d.Range($r.Start, $r.End - 1)\n    $rr.InsertXML($xmlFrag)\n}\n"}
